# backlog.xlsx - "primera implementacion de playerManager" edit
#
# Summary of the change being applied (per the commit/diff):
#  - Hoja2 is renamed to "Agregar jugadores" and becomes the active sheet.
#  - Hoja1!B6 ("Poder agregar jugadores") gets a hyperlink to the new sheet,
#    which also gives it the built-in Hyperlink style (underline + themed
#    color) and moves the sheet1 selection to B6.
#  - The new "Agregar jugadores" sheet is populated with three backlog
#    rows plus a "*" marker, and its B column is sized to fit the text.
#  - The new sheet's selection ends on C5.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename Hoja2 -> "Agregar jugadores"
$ws2.Name = "Agregar jugadores"

# 2. Add the new backlog rows on the "Agregar jugadores" sheet
$ws2.Range("B3").Value = "poner htmls en archivos separados"
$ws2.Range("B4").Value = "poder configurar por medio de un json"
$ws2.Range("B5").Value = "generar clase que permisista los jugadores"
$ws2.Range("C5").Value = "*"

# Size column B so the new text fits
$ws2.Columns.Item(2).AutoFit()

# 3. Turn Hoja1!B6 into a hyperlink pointing at the new sheet
$ws1.Hyperlinks.Add($ws1.Range("B6"), "", "'Agregar jugadores'!A1", "", "Poder agregar jugadores")

# 4. Update the selection on Hoja1 to B6
$ws1.Range("B6").Select()

# 5. Make "Agregar jugadores" the active sheet/tab, selected at C5
$ws2.Activate()
$ws2.Range("C5").Select()
